# Timesheet update: fix week-1 date labels (they had been stored as
# mis-parsed date serials instead of text like the other weeks) and record
# the work done on 24 September 2018 (Monday), plus a couple of related
# entries that were filled in at the same time (19 Sep "off"-adjacent notes
# for 20/21 Sep, and 22/23 Sep marked "off").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 1 (rows 20-26): replace numeric date serials with plain text
# dates "1.9.18".."7.9.18", matching the style used by the other weeks.
# Toggling the NumberFormat to Text and back keeps the original cell
# style (s="24") while preventing Excel from re-parsing the text back
# into a date serial.
$week1Dates = @{
    20 = "1.9.18"
    21 = "2.9.18"
    22 = "3.9.18"
    23 = "4.9.18"
    24 = "5.9.18"
    25 = "6.9.18"
    26 = "7.9.18"
}
foreach ($r in 20..26) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $week1Dates[$r]
    $cell.NumberFormat = "M/D/YYYY"
}

# --- Row 44 : Wednesday 19.9.18 (3rd week) - add time in/out and task
$ws.Range("G44").Value = 10/24
$ws.Range("H44").Value = 16/24
$ws.Range("J44").Value = "learnt about odoo "

# --- Rows 45-46 : Thursday 20.9.18 and Friday 21.9.18 (3rd week) - mark off
$ws.Range("J45").Value = "off"
$ws.Range("J46").Value = "off"

# --- Rows 50-51 : Saturday 22.9.18 and Sunday 23.9.18 (4th week) - mark off
$ws.Range("J50").Value = "off"
$ws.Range("J51").Value = "off"

# --- Row 52 : Monday 24.9.18 (4th week) - add time in/out and task
$ws.Range("G52").Value = 10/24
$ws.Range("H52").Value = 18/24
$ws.Range("J52").Value = "tried some basic commmands of linux CLI and started php"

# --- Restore view to top of sheet with the edited cell selected, mirroring
# the author's final cursor position after making the edit.
$ws.Range("J53").Select() | Out-Null
